$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 38.198408
$ws.Range("N2").Value = 114.595224
$ws.Range("O2").Value = 0.05104410684611114
$ws.Range("P2").Value = 0.05104410684611114
$ws.Range("Q2").Value = 0.167729209528
$ws.Range("R2").Value = 1.509562885752
$ws.Range("S2").Value = 0.05104410684611114
$ws.Range("T2").Value = 0.05104410684611114

# Row 3
$ws.Range("O3").Value = 0.2557395719837403
$ws.Range("P3").Value = 0.2557395719837403
$ws.Range("R3").Value = 7.563164293356001
$ws.Range("S3").Value = 0.2557395719837403
$ws.Range("T3").Value = 0.2557395719837403

# Row 4
$ws.Range("M4").Value = 61.421814
$ws.Range("N4").Value = 184.265442
$ws.Range("O4").Value = 0.08207728543288938
$ws.Range("P4").Value = 0.08207728543288938
$ws.Range("Q4").Value = 0.269703185274
$ws.Range("R4").Value = 2.427328667466
$ws.Range("S4").Value = 0.08207728543288938
$ws.Range("T4").Value = 0.08207728543288938

# Row 5
$ws.Range("M5").Value = 116.4573846666667
$ws.Range("N5").Value = 349.372154
$ws.Range("O5").Value = 0.1556207050813216
$ws.Range("P5").Value = 0.1556207050813217
$ws.Range("Q5").Value = 0.5113643760713332
$ws.Range("R5").Value = 4.602279384641999
$ws.Range("S5").Value = 0.1556207050813216
$ws.Range("T5").Value = 0.1556207050813217

# Row 6
$ws.Range("M6").Value = 126.0955403333333
$ws.Range("N6").Value = 378.286621
$ws.Range("O6").Value = 0.1685000650705857
$ws.Range("P6").Value = 0.1685000650705857
$ws.Range("Q6").Value = 0.5536855176036666
$ws.Range("R6").Value = 4.983169658433
$ws.Range("S6").Value = 0.1685000650705857
$ws.Range("T6").Value = 0.1685000650705857

# Row 7
$ws.Range("M7").Value = 214.7875923333334
$ws.Range("N7").Value = 644.3627770000001
$ws.Range("O7").Value = 0.2870182655853519
$ws.Range("P7").Value = 0.2870182655853519
$ws.Range("Q7").Value = 0.9431323179356668
$ws.Range("R7").Value = 8.488190861421002
$ws.Range("S7").Value = 0.2870182655853519
$ws.Range("T7").Value = 0.2870182655853519
